# Actualización 11 de Mayo - Mañana
#
# The "Rescatables" sheet gets a new student record inserted as the new
# row 2 (NC 20330051920030 - TZANAHUA GONZALEZ XIMENA -
# "LECTURA, EXPRESIÓN ORAL Y ESCRITA II" / 2AEM / 2 reprobadas), pushing
# the previously-existing record (NC 19330051920321 - DE JESUS DE LA CRUZ
# IGNACIO) down to row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Columns, per header row: A=NC, B=Paterno, C=Materno, D=Nombres,
# E=Nombre_Largo, F=Grupo, G=Reprobadas

# Push the existing data row (currently row 2) down to row 3 first so we
# don't clobber it.
$ws.Range("A3").Value = $ws.Range("A2").Value()
$ws.Range("B3").Value = $ws.Range("B2").Value()
$ws.Range("C3").Value = $ws.Range("C2").Value()
$ws.Range("D3").Value = $ws.Range("D2").Value()
$ws.Range("E3").Value = $ws.Range("E2").Value()
$ws.Range("F3").Value = $ws.Range("F2").Value()
$ws.Range("G3").Value = $ws.Range("G2").Value()

# Write the new student record into row 2.
$ws.Range("A2").Value = 20330051920030
$ws.Range("B2").Value = "TZANAHUA"
$ws.Range("C2").Value = "GONZALEZ"
$ws.Range("D2").Value = "XIMENA"
$ws.Range("E2").Value = "LECTURA, EXPRESIÓN ORAL Y ESCRITA II"
$ws.Range("F2").Value = "2AEM"
$ws.Range("G2").Value = 2
